$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, shifting existing rows 183:213 down to 184:214.
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with a new weekly price entry
# (same market/product as its neighbours, new date 2021-11-22 = serial 44522).
$ws.Cells.Item(183, 1).Value = 3
$ws.Cells.Item(183, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(183, 3).Value = "Coquimbo"
$ws.Cells.Item(183, 4).Value = 44522
$ws.Cells.Item(183, 5).Value = 5
$ws.Cells.Item(183, 6).Value = 100112039
$ws.Cells.Item(183, 7).Value = "Ciboulette"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 160
$ws.Cells.Item(183, 11).Value = 1500
$ws.Cells.Item(183, 12).Value = 1500
$ws.Cells.Item(183, 13).Value = 1500
$ws.Cells.Item(183, 14).Value = '$/docena de atados'
$ws.Cells.Item(183, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(183, 16).Value = 500
$ws.Cells.Item(183, 17).Value = 3
$ws.Cells.Item(183, 18).Value = "Hortaliza"
